# Apply updated cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.349.46"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.882.53"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "246.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.35%  "
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +8.56%  "
$ws.Range("D14").Value = "2.156.91"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "1.910.28"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "35.317.07"
$ws.Range("D19").Value = "0.0₃0827"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.16%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0599"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -12.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.855"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0735"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.98%  "
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "1.306.01"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0809"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "2.060.44"
$ws.Range("E51").Value = "  +0.13%  "
